# Automatische test-sync: 2025-07-27 19:11:50
#
# This script narrows the "Logs" sheet down from 4 test rows to just the
# first one (and fills in its reply / flag columns), and updates the
# "Dashboard" summary sheet + its chart to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: fill in row 2's answer + update its tracking columns,
#    then drop rows 3-5 (the other three test mails) entirely.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("E2").Value = "Geachte klant,`nDank u voor uw bericht. Om u zo goed mogelijk van dienst te zijn, zou u wat meer details kunnen geven over wat u precies geregeld wilt hebben? Zo kunnen wij u beter helpen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F2").Value = "2025-07-27 19:11:15"
$logs.Range("G2").Value = "Ja"
$logs.Range("H2").Value = "Nee"
$logs.Range("I2").Value = "Ja"
$logs.Range("J2").Value = "Nee"

# Writing the multi-line answer bumps the row to a custom auto-fit height;
# auto-fit it back down to the sheet default so row 2 keeps its plain height.
$logs.Rows.Item(2).AutoFit()

# Remove rows 5, 4, 3 (bottom-up so row indices don't shift underneath us).
$logs.Rows.Item(5).Delete()
$logs.Rows.Item(4).Delete()
$logs.Rows.Item(3).Delete()

# Row deletion doesn't auto-shrink the conditional formatting ranges, so
# pull each one back down to just row 2.
$condCols = @("D", "G", "H", "I", "J")
foreach ($col in $condCols) {
    $cellAddr = $col + "2"
    $fcs = $logs.Range($cellAddr).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($cellAddr))
    }
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: only one category remains ("Overig"), so its count
#    drops from 2 to 1 and the second summary row goes away.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 1
$dash.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# 3) Chart: category/value series now only cover the single remaining row.
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$series = $chart.SeriesCollection(1)
$series.Values = $dash.Range("B2")
$series.XValues = $dash.Range("A2")
